$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows for the "Extended Dog Talents" mabari abilities, right
# before the old row 97 (so old rows 97..473 shift down to 105..481).
$ws.Rows("97:104").Insert()

# Copy number/fill formatting from an existing "highlighted numeric id" row
# (A60, style index 6) onto the new A97:A104 cells so they render the same
# way (light fill, General number format) as the rest of the new-content
# blocks added elsewhere in the sheet.
$ws.Range("A60").Copy()
$ws.Range("A97:A104").PasteSpecial(-4122)

# Set the A-column (numeric id) cells first; order doesn't affect the
# shared-string table since these are stored as inline numbers.
$ws.Range("A97").Value = 6610096
$ws.Range("A98").Value = 6610097
$ws.Range("A99").Value = 6610098
$ws.Range("A100").Value = 6610099
$ws.Range("A101").Value = 6610100
$ws.Range("A102").Value = 6610101
$ws.Range("A103").Value = 6610102
$ws.Range("A104").Value = 6610103

# Set the B-column (text) cells in the exact order the strings were first
# authored, so new entries land in the shared-string table in that order:
# Endurance, Endurance-desc, Bond, Frighten, Ferocious bite, Bond-desc,
# Frighten-desc, Ferocious-bite-desc.
$ws.Range("B97").Value = "Endurance"
$ws.Range("B98").Value = "The mabari has undergone endurance and survival training, gaining a bonus to stamina."
$ws.Range("B99").Value = "Bond"
$ws.Range("B101").Value = "Frighten"
$ws.Range("B103").Value = "Ferocious bite"
$ws.Range("B100").Value = "The mabari has formed a special bond with its master, gaining a bonus to mental resistance."
$ws.Range("B102").Value = "The mabari lets out a frightening howl, forcing nearby enemies to cower in fear unless they pass a mental resistance check."
$ws.Range("B104").Value = "The mabari will jump on its target and bite its neck. If the target is a living creature, it will die instantly if of normal or lesser rank unless it passes a physical resistance check. Lieutenant-ranked enemies will suffer a critical hit, while boss-ranked enemies will take only standard damage. Additionally, all targets who can bleed will take additional damage over time.`nNon-bleeding creatures will only take normal damage."

# Row 104 (the long "Ferocious bite" description) is taller than the rest.
$ws.Rows("104:104").RowHeight = 45

# Update the "strings" defined name so it spans the new, larger data range.
$wb.Names.Item("strings").RefersTo = "=Sheet1!`$A`$1:`$B`$473"

# Move selection to mirror where the author was last working.
$ws.Range("B93").Select()
